$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cells D2:E2 and D3:E3
$ws.Range("D2").Value = -0.0692
$ws.Range("E2").Value = -0.0469
$ws.Range("D3").Value = -0.0692
$ws.Range("E3").Value = -0.0469

# Update company name (legal entity suffix Q.S.C. -> Q.P.S.C.)
$ws.Range("B3").Value = "Mesaieed Petrochemical Holding Company Q.P.S.C. (DSM:MPHC)"

# Updated capital-structure derived metrics for row 2
$ws.Range("I2").Value = 0.979947689625109
$ws.Range("J2").Value = 0.979947689625109
$ws.Range("K2").Value = 280.7
$ws.Range("L2").Value = 1.223626852659111
$ws.Range("M2").Value = 241.5
$ws.Range("N2").Value = 0.03418694525841933
$ws.Range("O2").Value = 0.8603491271820449
$ws.Range("P2").Value = 241.5
$ws.Range("Q2").Value = 0.03418694525841933
$ws.Range("R2").Value = 0.8603491271820449
$ws.Range("U2").Value = 75.8
$ws.Range("V2").Value = 0.0107303124247958
$ws.Range("W2").Value = 0.07115335868187579
$ws.Range("X2").Value = 0.05870255866550585
$ws.Range("Y2").Value = 0.01245080001636995
$ws.Range("Z2").Value = 0.05817255999837705
$ws.Range("AA2").Value = 0.05700606576998762
$ws.Range("AB2").Value = 0.05870255866550585
$ws.Range("AC2").Value = -0.001696492895518223
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = -75.8
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = -0.01084670091438547
$ws.Range("AK2").Value = -0.01825406381697772

# Updated capital-structure derived metrics for row 3
$ws.Range("I3").Value = 0.979947689625109
$ws.Range("J3").Value = 0.979947689625109
$ws.Range("K3").Value = 280.7
$ws.Range("L3").Value = 1.223626852659111
$ws.Range("M3").Value = 241.5
$ws.Range("N3").Value = 0.03418694525841933
$ws.Range("O3").Value = 0.8603491271820449
$ws.Range("P3").Value = 241.5
$ws.Range("Q3").Value = 0.03418694525841933
$ws.Range("R3").Value = 0.8603491271820449
$ws.Range("U3").Value = 75.8
$ws.Range("V3").Value = 0.0107303124247958
$ws.Range("W3").Value = 0.07115335868187579
$ws.Range("X3").Value = 0.05870255866550585
$ws.Range("Y3").Value = 0.01245080001636995
$ws.Range("Z3").Value = 0.05817255999837705
$ws.Range("AA3").Value = 0.05700606576998762
$ws.Range("AB3").Value = 0.05870255866550585
$ws.Range("AC3").Value = -0.001696492895518223
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = -75.8
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = -0.01084670091438547
$ws.Range("AK3").Value = -0.01825406381697772

# Remove debt_ebitda (AN) and net_debt_ebitda (AP) values for rows 2 and 3
$ws.Range("AN2").ClearContents()
$ws.Range("AN3").ClearContents()
$ws.Range("AP2").ClearContents()
$ws.Range("AP3").ClearContents()
